$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Unprotect()
$ws.Range("A1").Value = "TEST VALUE"
$ws.Protect($null, $true, $true, $true, $null, $null, $null, $null, $null, $true, $null, $null, $true)
